$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 2, 5, and 10 as "fertig" (done) in the Status column (D)
$ws.Range("D2").Value = "fertig"
$ws.Range("D5").Value = "fertig"
$ws.Range("D10").Value = "fertig"

# Update the current selection to D17 (as last edited/viewed cell)
$ws.Range("D17").Select()
